$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add new row 21, a duplicate of the original (pre-edit) row 20 ---
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44316
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = "Frutos de pepita"
$ws.Range("I21").Value = 100104003
$ws.Range("J21").Value = "Membrillo"
$ws.Range("K21").Value = "Champion"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9500
$ws.Range("Q21").Value = "$/caja 18 kilos granel"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 528
$ws.Range("T21").Value = 18

# --- Step 2: update row 20 with the new values from the diff ---
$ws.Range("D20").Value = 45034
$ws.Range("M20").Value = 220
$ws.Range("N20").Value = 8500
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 8727
$ws.Range("S20").Value = 485
